$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 183
$ws.Range("I18").Value = 183
$ws.Range("J18").Value = 0
$ws.Range("K18").Value = 183
$ws.Range("L18").Value = 0
$ws.Range("M18").Value = 101
$ws.Range("N18").ClearContents()

$ws.Range("H62").Value = 4470.8887
$ws.Range("I62").Value = 2868.4
$ws.Range("J62").Value = 6474
$ws.Range("K62").Value = 2868.4
$ws.Range("L62").Value = 6474
$ws.Range("M62").Value = -2244.4
$ws.Range("N62").ClearContents()

$ws.Range("H64").Value = 100002500
$ws.Range("I64").Value = 200000000
$ws.Range("J64").Value = 4999
$ws.Range("K64").Value = 200000000
$ws.Range("L64").Value = 4999
$ws.Range("M64").Value = -199999752
$ws.Range("N64").ClearContents()

$ws.Range("H65").Value = 4470.8887
$ws.Range("I65").Value = 2868.4
$ws.Range("J65").Value = 6474
$ws.Range("K65").Value = 14342
$ws.Range("L65").Value = 32370
$ws.Range("M65").Value = -11222
$ws.Range("N65").ClearContents()

$ws.Range("H67").Value = 100002500
$ws.Range("I67").Value = 200000000
$ws.Range("J67").Value = 4999
$ws.Range("K67").Value = 200000000
$ws.Range("L67").Value = 4999
$ws.Range("M67").Value = -199999142
$ws.Range("N67").ClearContents()

$ws.Range("H70").Value = 11199.857
$ws.Range("I70").Value = 1932.6666
$ws.Range("J70").Value = 13727.272
$ws.Range("K70").Value = 5797.9998
$ws.Range("L70").Value = 41181.81600000001
$ws.Range("M70").Value = -5527.9998
$ws.Range("N70").ClearContents()

$ws.Range("H73").Value = 11199.857
$ws.Range("I73").Value = 1932.6666
$ws.Range("J73").Value = 13727.272
$ws.Range("K73").Value = 5797.9998
$ws.Range("L73").Value = 41181.81600000001
$ws.Range("M73").Value = -4861.9998
$ws.Range("N73").ClearContents()

$ws.Range("H76").Value = 4563.9287
$ws.Range("I76").Value = 4098.8887
$ws.Range("J76").Value = 5401
$ws.Range("K76").Value = 4098.8887
$ws.Range("L76").Value = 5401
$ws.Range("M76").Value = -3783.8887
$ws.Range("N76").ClearContents()

$ws.Range("H79").Value = 4563.9287
$ws.Range("I79").Value = 4098.8887
$ws.Range("J79").Value = 5401
$ws.Range("K79").Value = 4098.8887
$ws.Range("L79").Value = 5401
$ws.Range("M79").Value = -3006.8887
$ws.Range("N79").ClearContents()

$ws.Range("H106").Value = 2357.3
$ws.Range("I106").Value = 2218.6
$ws.Range("J106").Value = 2496
$ws.Range("K106").Value = 2218.6
$ws.Range("L106").Value = 2496
$ws.Range("M106").Value = -1587.6
$ws.Range("N106").Value = -3758

$ws.Range("H118").Value = 2185.4
$ws.Range("I118").Value = 2185.4
$ws.Range("J118").Value = 0
$ws.Range("K118").Value = 6556.200000000001
$ws.Range("L118").Value = 0
$ws.Range("M118").Value = -4899.200000000001
$ws.Range("N118").ClearContents()

$ws.Range("H137").Value = 2295.3
$ws.Range("I137").Value = 1478.6
$ws.Range("J137").Value = 4745.4
$ws.Range("K137").Value = 4435.799999999999
$ws.Range("L137").Value = 14236.2
$ws.Range("M137").Value = -1885.799999999999
$ws.Range("N137").ClearContents()

$ws.Range("H138").Value = 1280.3654
$ws.Range("I138").Value = 857.6667
$ws.Range("J138").Value = 3997.7144
$ws.Range("K138").Value = 2573.0001
$ws.Range("L138").Value = 11993.1432
$ws.Range("M138").Value = 2566.9999
$ws.Range("N138").Value = -22273.1432

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 547274.1
$ws.Range("I2").Value = 641822.25
$ws.Range("J2").Value = 3622.25
$ws.Range("K2").Value = 641822.25
$ws.Range("L2").Value = 3622.25
$ws.Range("M2").Value = -641709.25
$ws.Range("N2").ClearContents()

$ws.Range("H43").Value = 0
$ws.Range("I43").Value = 0
$ws.Range("J43").Value = 0
$ws.Range("K43").Value = 0
$ws.Range("L43").Value = 0
$ws.Range("N43").ClearContents()

$ws.Range("H76").Value = 24977.6
$ws.Range("I76").Value = 0
$ws.Range("J76").Value = 24977.6
$ws.Range("K76").Value = 0
$ws.Range("L76").Value = 24977.6
$ws.Range("N76").Value = -25653.6

$ws.Range("H79").Value = 24977.6
$ws.Range("I79").Value = 0
$ws.Range("J79").Value = 24977.6
$ws.Range("K79").Value = 0
$ws.Range("L79").Value = 24977.6
$ws.Range("N79").Value = -27317.6

$ws.Range("H97").Value = 467.08334
$ws.Range("I97").Value = 491.72726
$ws.Range("J97").Value = 196
$ws.Range("K97").Value = 491.72726
$ws.Range("L97").Value = 196
$ws.Range("M97").Value = 4.272739999999999
$ws.Range("N97").ClearContents()

$ws.Range("H116").Value = 547274.1
$ws.Range("I116").Value = 641822.25
$ws.Range("J116").Value = 3622.25
$ws.Range("K116").Value = 641822.25
$ws.Range("L116").Value = 3622.25
$ws.Range("M116").Value = -639528.25
$ws.Range("N116").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 547274.1
$ws.Range("I3").Value = 641822.25
$ws.Range("J3").Value = 3622.25
$ws.Range("K3").Value = 641822.25
$ws.Range("L3").Value = 3622.25
$ws.Range("M3").Value = -641708.25
$ws.Range("N3").ClearContents()

$ws.Range("H7").Value = 2500499
$ws.Range("I7").Value = 3333666
$ws.Range("J7").Value = 998
$ws.Range("K7").Value = 3333666
$ws.Range("L7").Value = 998
$ws.Range("M7").Value = -3333553
$ws.Range("N7").ClearContents()

$ws.Range("H107").Value = 78374.62
$ws.Range("I107").Value = 1492.8182
$ws.Range("J107").Value = 501224.5
$ws.Range("K107").Value = 1492.8182
$ws.Range("L107").Value = 501224.5
$ws.Range("M107").Value = 427.1818000000001
$ws.Range("N107").Value = -505064.5

$ws.Range("H134").Value = 27107666
$ws.Range("I134").Value = 30296426
$ws.Range("J134").Value = 3199.5
$ws.Range("K134").Value = 90889278
$ws.Range("L134").Value = 9598.5
$ws.Range("M134").Value = -90886743
$ws.Range("N134").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5696.6787
$ws.Range("I31").Value = 4856.8
$ws.Range("J31").Value = 6163.278
$ws.Range("K31").Value = 4856.8
$ws.Range("L31").Value = 6163.278
$ws.Range("M31").Value = -4561.8
$ws.Range("N31").Value = -6753.278

$ws.Range("H34").Value = 5696.6787
$ws.Range("I34").Value = 4856.8
$ws.Range("J34").Value = 6163.278
$ws.Range("K34").Value = 4856.8
$ws.Range("L34").Value = 6163.278
$ws.Range("M34").Value = -4654.8
$ws.Range("N34").Value = -6567.278

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 46397.363
$ws.Range("I5").Value = 84049.5
$ws.Range("J5").Value = 1214.8
$ws.Range("K5").Value = 252148.5
$ws.Range("L5").Value = 3644.4
$ws.Range("M5").Value = -252036.5
$ws.Range("N5").Value = -3868.4

$ws.Range("H120").Value = 4722
$ws.Range("I120").Value = 5000
$ws.Range("J120").Value = 4444
$ws.Range("K120").Value = 15000
$ws.Range("L120").Value = 13332
$ws.Range("M120").Value = -10162
$ws.Range("N120").Value = -23008

$ws.Range("H135").Value = 46397.363
$ws.Range("I135").Value = 84049.5
$ws.Range("J135").Value = 1214.8
$ws.Range("K135").Value = 756445.5
$ws.Range("L135").Value = 10933.2
$ws.Range("M135").Value = -753910.5
$ws.Range("N135").Value = -16003.2

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 8748.5
$ws.Range("I43").Value = 6664.6665
$ws.Range("J43").Value = 15000
$ws.Range("K43").Value = 6664.6665
$ws.Range("L43").Value = 15000
$ws.Range("M43").Value = -6513.6665
$ws.Range("N43").Value = -15302

$ws.Range("H46").Value = 3999.5
$ws.Range("I46").Value = 3999.5
$ws.Range("J46").Value = 0
$ws.Range("K46").Value = 3999.5
$ws.Range("L46").Value = 0
$ws.Range("M46").Value = -3843.5

$ws.Range("H57").Value = 9999
$ws.Range("I57").Value = 9999
$ws.Range("J57").Value = 0
$ws.Range("K57").Value = 9999
$ws.Range("L57").Value = 0
$ws.Range("M57").Value = -9179

$ws.Range("H80").Value = 2457
$ws.Range("I80").Value = 2457
$ws.Range("J80").Value = 0
$ws.Range("K80").Value = 2457
$ws.Range("L80").Value = 0
$ws.Range("M80").Value = -1459

$ws.Range("H83").Value = 2457
$ws.Range("I83").Value = 2457
$ws.Range("J83").Value = 0
$ws.Range("K83").Value = 12285
$ws.Range("L83").Value = 0
$ws.Range("M83").Value = -7293

$ws.Range("H132").Value = 5210912.5
$ws.Range("I132").Value = 5683972
$ws.Range("J132").Value = 7257
$ws.Range("K132").Value = 17051916
$ws.Range("L132").Value = 21771
$ws.Range("M132").Value = -17049386
$ws.Range("N132").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1954.3043
$ws.Range("I46").Value = 1895.8334
$ws.Range("J46").Value = 2018.091
$ws.Range("K46").Value = 1895.8334
$ws.Range("L46").Value = 2018.091
$ws.Range("M46").Value = -1707.8334
$ws.Range("N46").Value = -2394.091

$ws.Range("H55").Value = 473.13333
$ws.Range("I55").Value = 217.3
$ws.Range("J55").Value = 984.8
$ws.Range("K55").Value = 217.3
$ws.Range("L55").Value = 984.8
$ws.Range("M55").Value = -44.30000000000001
$ws.Range("N55").Value = -1330.8

$ws.Range("H128").Value = 25000
$ws.Range("I128").Value = 0
$ws.Range("J128").Value = 25000
$ws.Range("K128").Value = 0
$ws.Range("L128").Value = 25000
$ws.Range("N128").Value = -34960

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 6017.4
$ws.Range("I62").Value = 3044.75
$ws.Range("J62").Value = 7999.1665
$ws.Range("K62").Value = 3044.75
$ws.Range("L62").Value = 7999.1665
$ws.Range("M62").Value = -2420.75
$ws.Range("N62").ClearContents()

$ws.Range("H65").Value = 6017.4
$ws.Range("I65").Value = 3044.75
$ws.Range("J65").Value = 7999.1665
$ws.Range("K65").Value = 15223.75
$ws.Range("L65").Value = 39995.8325
$ws.Range("M65").Value = -12103.75
$ws.Range("N65").ClearContents()

$ws.Range("H122").Value = 2249.5
$ws.Range("I122").Value = 2174.75
$ws.Range("J122").Value = 2399
$ws.Range("K122").Value = 6524.25
$ws.Range("L122").Value = 7197
$ws.Range("M122").Value = -4074.25
$ws.Range("N122").Value = -12097

$ws.Range("H132").Value = 15627295
$ws.Range("I132").Value = 19232926
$ws.Range("J132").Value = 2893
$ws.Range("K132").Value = 57698778
$ws.Range("L132").Value = 8679
$ws.Range("M132").Value = -57696248
$ws.Range("N132").Value = -13184.8568
